# Refactor printer output: reword the instructions cell on the English sheet
# and rename the Hungarian "Szám" header to "Válasz" (Answer).
#
# Order matters for the shared-string table: the English instructions text
# must be written before the Hungarian "Válasz" label so the two brand-new
# shared strings get appended in the same order as the target workbook.

$wb   = $excel.ActiveWorkbook
$wsHU = $wb.Worksheets.Item("MintaHU")
$wsEN = $wb.Worksheets.Item("MintaEN")

# MintaEN: reword the instructions cell (A4)
$wsEN.Range("A4").Value = "Write the number next to the name of each building, that is shown near its model in the exhibition!"

# MintaHU: rename the "Szám" column header (C5) to "Válasz"
$wsHU.Range("C5").Value = "Válasz"

# MintaEN: move the selection/active cell to A1:C1 (was A3:C3), without
# leaving MintaEN as the active tab - restore MintaHU as the active sheet
# afterwards, matching the original workbook state.
$wsEN.Range("A1:C1").Select()
$wsHU.Activate()
